$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to
# text so Excel keeps them as strings (matching the original inline-string
# cell type) instead of silently converting to a numeric value.

$ws.Range("D2").Value = "25.778.21"
$ws.Range("E2").Value = "  -0.85%  "
$ws.Range("D3").Value = "1.598.31"
$ws.Range("E3").Value = "  -2.41%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.71%  "
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.481"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.69%  "
$ws.Range("E9").Value = "  -2.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.78"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.00%  "
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("D12").Value = "1.819.67"
$ws.Range("E12").Value = "  -2.46%  "
$ws.Range("D13").Value = "1.601.14"
$ws.Range("E13").Value = "  -2.17%  "
$ws.Range("E14").Value = "  -4.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.507"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.33%  "
$ws.Range("D16").Value = "25.760.70"
$ws.Range("E16").Value = "  -0.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.39"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.42%  "
$ws.Range("D18").Value = "0.0₃0713"
$ws.Range("E18").Value = "  -4.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.01"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "188.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.60%  "
$ws.Range("E21").Value = "  -1.89%  "
$ws.Range("E22").Value = "  -4.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.09%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("E25").Value = "  -3.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.78"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.16%  "
$ws.Range("E27").Value = "  -4.73%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.50"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.32%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "14.93"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.19%  "
$ws.Range("E30").Value = "  -4.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0467"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.35%  "
$ws.Range("E32").Value = "  -3.04%  "
$ws.Range("E33").Value = "  -5.06%  "
$ws.Range("E34").Value = "  -0.45%  "
$ws.Range("E35").Value = "  -1.91%  "
$ws.Range("D36").Value = "1.094.12"
$ws.Range("E36").Value = "  -3.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.36"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.16%  "
$ws.Range("E38").Value = "  -0.41%  "
$ws.Range("E39").Value = "  -2.46%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.788"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -9.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.495"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "95.25"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.30%  "
$ws.Range("D43").Value = "1.735.13"
$ws.Range("E43").Value = "  -2.25%  "
$ws.Range("E44").Value = "  -3.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.740"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.01%  "
$ws.Range("E46").Value = "  -1.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "53.16"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.78%  "
$ws.Range("E48").Value = "  -3.25%  "
$ws.Range("E49").Value = "  -1.02%  "
$ws.Range("E50").Value = "  -5.60%  "
$ws.Range("E51").Value = "  -0.16%  "
